$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.296300000000006
$ws.Range("A12").Value = -21.5677
$ws.Range("B14").Value = 6.898399999999998
$ws.Range("B26").Value = 3.903600000000003
$ws.Range("A27").Value = -21.72609999999999
$ws.Range("B31").Value = 5.653700000000004
$ws.Range("A32").Value = -21.28640000000001
$ws.Range("B35").Value = 9.343400000000006
$ws.Range("A36").Value = -19.9758
$ws.Range("B37").Value = 8.864600000000003
$ws.Range("A38").Value = -19.4298
$ws.Range("B45").Value = 7.184500000000001
$ws.Range("A46").Value = -21.5711
$ws.Range("B52").Value = 5.244400000000002
$ws.Range("A54").Value = -21.68129999999999
$ws.Range("A55").Value = -22.4205
$ws.Range("A56").Value = -22.05919999999999
$ws.Range("B57").Value = 4.845999999999996
$ws.Range("A67").Value = -21.44239999999998
$ws.Range("A69").Value = -21.58419999999998
$ws.Range("A72").Value = -22.04060000000001
$ws.Range("B81").Value = 6.707600000000001
$ws.Range("A83").Value = -21.67679999999999
$ws.Range("B83").Value = 5.317400000000001
$ws.Range("A86").Value = -21.9531
$ws.Range("A91").Value = -21.3913
$ws.Range("A93").Value = -21.3123
$ws.Range("A99").Value = -20.03759999999999
$ws.Range("B100").Value = 5.036299999999998
$ws.Range("B102").Value = 8.390700000000002
